$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.892066666666667
$ws.Range("H2").Value = 23.6762
$ws.Range("I2").Value = 0.1739002798877711
$ws.Range("J2").Value = 0.1739002798877711
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 29.9519109074889
$ws.Range("R2").Value = 269.5671981674
$ws.Range("S2").Value = 0.001815890643810962
$ws.Range("T2").Value = 0.001815890643810962
$ws.Range("G3").Value = 7.892066666666667
$ws.Range("H3").Value = 23.6762
$ws.Range("I3").Value = 0.1739002798877711
$ws.Range("J3").Value = 0.1739002798877711
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 1920.742763370578
$ws.Range("R3").Value = 17286.6848703352
$ws.Range("S3").Value = 0.1164486240609167
$ws.Range("T3").Value = 0.1164486240609167
$ws.Range("G4").Value = 7.892066666666667
$ws.Range("H4").Value = 23.6762
$ws.Range("I4").Value = 0.1739002798877711
$ws.Range("J4").Value = 0.1739002798877711
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 235.1944066900667
$ws.Range("R4").Value = 2116.7496602106
$ws.Range("S4").Value = 0.01425910099373251
$ws.Range("T4").Value = 0.01425910099373251
$ws.Range("G5").Value = 7.892066666666667
$ws.Range("H5").Value = 23.6762
$ws.Range("I5").Value = 0.1739002798877711
$ws.Range("J5").Value = 0.1739002798877711
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 682.4806128448445
$ws.Range("R5").Value = 6142.325515603601
$ws.Range("S5").Value = 0.04137666418931085
$ws.Range("T5").Value = 0.04137666418931085
$ws.Range("I6").Value = 0.3815924715300191
$ws.Range("J6").Value = 0.3815924715300191
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 65.72400986135135
$ws.Range("R6").Value = 591.5160887521621
$ws.Range("S6").Value = 0.003984641078480464
$ws.Range("T6").Value = 0.003984641078480465
$ws.Range("I7").Value = 0.3815924715300191
$ws.Range("J7").Value = 0.3815924715300191
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("R7").Value = 37932.47951348279
$ws.Range("S7").Value = 0.2555252831700593
$ws.Range("T7").Value = 0.2555252831700593
$ws.Range("I8").Value = 0.3815924715300191
$ws.Range("J8").Value = 0.3815924715300191
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 516.091262169442
$ws.Range("R8").Value = 4644.821359524978
$ws.Range("S8").Value = 0.03128899846225705
$ws.Range("T8").Value = 0.03128899846225704
$ws.Range("I9").Value = 0.3815924715300191
$ws.Range("J9").Value = 0.3815924715300191
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 1497.579325317119
$ws.Range("R9").Value = 13478.21392785407
$ws.Range("S9").Value = 0.09079354881922228
$ws.Range("T9").Value = 0.09079354881922229
$ws.Range("G10").Value = 7.716272666666666
$ws.Range("H10").Value = 23.148818
$ws.Range("I10").Value = 0.1700266904854272
$ws.Range("J10").Value = 0.1700266904854272
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 29.28473886644289
$ws.Range("R10").Value = 263.562649797986
$ws.Range("S10").Value = 0.001775442090431859
$ws.Range("T10").Value = 0.001775442090431859
$ws.Range("G11").Value = 7.716272666666666
$ws.Range("H11").Value = 23.148818
$ws.Range("I11").Value = 0.1700266904854272
$ws.Range("J11").Value = 0.1700266904854272
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 1877.95865274337
$ws.Range("R11").Value = 16901.62787469033
$ws.Range("S11").Value = 0.1138547572979018
$ws.Range("T11").Value = 0.1138547572979018
$ws.Range("G12").Value = 7.716272666666666
$ws.Range("H12").Value = 23.148818
$ws.Range("I12").Value = 0.1700266904854272
$ws.Range("J12").Value = 0.1700266904854272
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 229.9555044764926
$ws.Range("R12").Value = 2069.599540288434
$ws.Range("S12").Value = 0.01394148274417064
$ws.Range("T12").Value = 0.01394148274417064
$ws.Range("G13").Value = 7.716272666666666
$ws.Range("H13").Value = 23.148818
$ws.Range("I13").Value = 0.1700266904854272
$ws.Range("J13").Value = 0.1700266904854272
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 667.2785115548004
$ws.Range("R13").Value = 6005.506603993204
$ws.Range("S13").Value = 0.04045500835292294
$ws.Range("T13").Value = 0.04045500835292295
$ws.Range("G14").Value = 12.45667266666667
$ws.Range("H14").Value = 37.370018
$ws.Range("I14").Value = 0.2744805580967825
$ws.Range("J14").Value = 0.2744805580967826
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 47.27546860337623
$ws.Range("R14").Value = 425.4792174303861
$ws.Range("S14").Value = 0.002866163744403546
$ws.Range("T14").Value = 0.002866163744403547
$ws.Range("G15").Value = 12.45667266666667
$ws.Range("H15").Value = 37.370018
$ws.Range("I15").Value = 0.2744805580967825
$ws.Range("J15").Value = 0.2744805580967826
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 3031.660132982837
$ws.Range("R15").Value = 27284.94119684553
$ws.Range("S15").Value = 0.1838000683062185
$ws.Range("T15").Value = 0.1838000683062186
$ws.Range("G16").Value = 12.45667266666667
$ws.Range("H16").Value = 37.370018
$ws.Range("I16").Value = 0.2744805580967825
$ws.Range("J16").Value = 0.2744805580967826
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 371.2259235648927
$ws.Range("R16").Value = 3341.033312084034
$ws.Range("S16").Value = 0.02250626624203216
$ws.Range("T16").Value = 0.02250626624203216
$ws.Range("G17").Value = 12.45667266666667
$ws.Range("H17").Value = 37.370018
$ws.Range("I17").Value = 0.2744805580967825
$ws.Range("J17").Value = 0.2744805580967826
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 1077.213099511867
$ws.Range("R17").Value = 9694.917895606804
$ws.Range("S17").Value = 0.06530805980412827
$ws.Range("T17").Value = 0.06530805980412828
